# Apply the Thu Jan 5 11:06:18 UTC 2023 GitHub Actions symbol-list refresh
# (cryptos.xlsx): updates Price (D), Volume 1h % (E) and Hora (G) for rows
# 2-51, plus a Coin/Link swap between rows 41 and 42 (KickToken <-> BKEXToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E/G hold numeric-looking values (prices, percentages, the hour) but must
# stay stored as text, matching the workbook's existing inlineStr cells.
# A leading apostrophe forces Excel to keep the literal text; re-applying the
# "Normal" cell style afterwards clears the quote-prefix formatting so the
# cell's style stays unchanged (matches the diff, which touches values only).
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "256.77"
Set-TextValue "E2" "0.66%"
Set-TextValue "G2" "11"

Set-TextValue "D3" "26.88"
Set-TextValue "E3" "-4.38%"
Set-TextValue "G3" "11"

Set-TextValue "D4" "4.744"
Set-TextValue "E4" "-9.43%"
Set-TextValue "G4" "11"

Set-TextValue "D5" "0.05930"
Set-TextValue "E5" "1.13%"
Set-TextValue "G5" "11"

Set-TextValue "D6" "6.665"
Set-TextValue "E6" "-0.72%"
Set-TextValue "G6" "11"

Set-TextValue "D7" "0.8667"
Set-TextValue "E7" "0.25%"
Set-TextValue "G7" "11"

Set-TextValue "D8" "0.9471"
Set-TextValue "E8" "-3.49%"
Set-TextValue "G8" "11"

Set-TextValue "D9" "0.1400"
Set-TextValue "E9" "-0.64%"
Set-TextValue "G9" "11"

Set-TextValue "D10" "0.03785"
Set-TextValue "E10" "8.55%"
Set-TextValue "G10" "11"

Set-TextValue "D11" "0.07145"
Set-TextValue "E11" "-0.25%"
Set-TextValue "G11" "11"

Set-TextValue "D12" "0.03167"
Set-TextValue "E12" "-0.60%"
Set-TextValue "G12" "11"

Set-TextValue "D13" "0.09243"
Set-TextValue "E13" "0.18%"
Set-TextValue "G13" "11"

Set-TextValue "D14" "0.001543"
Set-TextValue "E14" "-0.25%"
Set-TextValue "G14" "11"

Set-TextValue "D15" "0.0006067"
Set-TextValue "E15" "-0.42%"
Set-TextValue "G15" "11"

Set-TextValue "D16" "0.006050"
Set-TextValue "E16" "4.48%"
Set-TextValue "G16" "11"

Set-TextValue "D17" "3.496"
Set-TextValue "E17" "-0.41%"
Set-TextValue "G17" "11"

Set-TextValue "D18" "3.202"
Set-TextValue "E18" "-0.53%"
Set-TextValue "G18" "11"

Set-TextValue "D19" "2.219"
Set-TextValue "E19" "0.79%"
Set-TextValue "G19" "11"

Set-TextValue "D20" "0.3125"
Set-TextValue "E20" "-1.77%"
Set-TextValue "G20" "11"

Set-TextValue "E21" "0.37%"
Set-TextValue "G21" "11"

Set-TextValue "D22" "3.811"
Set-TextValue "E22" "7.57%"
Set-TextValue "G22" "11"

Set-TextValue "D23" "0.04216"
Set-TextValue "E23" "1.45%"
Set-TextValue "G23" "11"

Set-TextValue "G24" "11"

Set-TextValue "D25" "0.001224"
Set-TextValue "E25" "0.08%"
Set-TextValue "G25" "11"

Set-TextValue "E26" "-10.73%"
Set-TextValue "G26" "11"

Set-TextValue "E27" "-0.07%"
Set-TextValue "G27" "11"

Set-TextValue "E28" "1.87%"
Set-TextValue "G28" "11"

Set-TextValue "G29" "11"

Set-TextValue "G30" "11"

Set-TextValue "G31" "11"

Set-TextValue "G32" "11"

Set-TextValue "G33" "11"

Set-TextValue "G34" "11"

Set-TextValue "G35" "11"

Set-TextValue "G36" "11"

Set-TextValue "G37" "11"

Set-TextValue "G38" "11"

Set-TextValue "G39" "11"

Set-TextValue "D40" "0.03827"
Set-TextValue "E40" "0.47%"
Set-TextValue "G40" "11"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1102"
Set-TextValue "E41" "0.09%"
Set-TextValue "G41" "11"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.003968"
Set-TextValue "E42" "-31.35%"
Set-TextValue "G42" "11"

Set-TextValue "E43" "-4.38%"
Set-TextValue "G43" "11"

Set-TextValue "D44" "0.01142"
Set-TextValue "E44" "20.30%"
Set-TextValue "G44" "11"

Set-TextValue "E45" "5.15%"
Set-TextValue "G45" "11"

Set-TextValue "E46" "-0.10%"
Set-TextValue "G46" "11"

Set-TextValue "E47" "-26.29%"
Set-TextValue "G47" "11"

Set-TextValue "D48" "0.002431"
Set-TextValue "E48" "13.80%"
Set-TextValue "G48" "11"

Set-TextValue "E49" "-0.10%"
Set-TextValue "G49" "11"

Set-TextValue "E50" "-0.10%"
Set-TextValue "G50" "11"

Set-TextValue "G51" "11"

